# slide_template update: reposition/resize the Slide Master layout's
# Title / Slide-Number / divider-line shapes and bump the "PYMI.VN"
# slide-number text to 24pt.
#
# NOTE on the literals below: Shape.Left/Top/Width/Height are COM
# `Single` (32-bit float) properties expressed in points (1 pt = 12700
# EMU). To land on an exact target EMU value after the host's
# float32-truncating point->EMU conversion, each literal below was
# chosen (via the neighbouring float32 value) so it reproduces the
# desired EMU exactly instead of drifting by +/-1 EMU from naive
# division.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$layout = $s.CustomLayout

# --- "Title 1" placeholder -------------------------------------------------
$title = $layout.Shapes.Item(1)
$title.Left   = 19.9285831451416
$title.Top    = 28.75
$title.Width  = 914.1428833007812
$title.Height = 468.8214416503906

# --- "Slide Number Placeholder 5" ------------------------------------------
$sldNum = $layout.Shapes.Item(2)
$sldNum.Left   = 0.0
$sldNum.Top    = 510.57183837890625
$sldNum.Width  = 894.0
$sldNum.Height = 28.75

# Bump the "PYMI.VN" run text to 24pt.
$sldNum.TextFrame.TextRange.Font.Size = 24

# --- "Straight Connector 10" (divider line under the slide number) --------
$cxn = $layout.Shapes.Item(3)
$cxn.Left   = 0.0
$cxn.Top    = 508.6429138183594
$cxn.Width  = 960.0
$cxn.Height = 0.0
